# mk3-shield-PCBWay-BOM.xlsx update:
#  - PCB size text 60x52mm -> 60x60mm
#  - Row 11 item number becomes "5**" (do-not-order marker)
#  - Row 12 filled in with new JP4/6 SamTec connector info
#  - Row 14 filled in with "**" / "DO NOT ORDER, INCLUDE AND  INSTALL!!!" note
#  - New supporting rows 27-35 (blank spacer block + extra note row) added
#  - View zoom bumped to 110% on all three sheets, selection moved to G2 on Sheet1

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Title text update ---
$ws1.Range("D2").Value = "MK3-SHIELD   BOM  (Bill of Materials)  PCB: 60x60mm"

# --- Row 11: mark item 5 as do-not-order ("5**") ---
$ws1.Range("A11").Value = "5**"

# --- Row 12: new JP4/6 SamTec stacking connector line ---
$ws1.Range("A12").Value = "6**"
$ws1.Range("B12").Value = "JP4/6"
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = "SamTec"
$ws1.Range("E12").Value = "ESQ-122-13-L-T"
$ws1.Range("F12").Value = "3row 44pos in two 2x 3x22 ESQ-122-13-L-T stacking board inter connect"
$ws1.Range("H12").Value = "Through Hole"
$ws1.Range("I12").Value = "do not populate (bottom side mount!!)"

# --- Row 14: legend note for the "**" marker ---
$ws1.Range("A14").Value = "**"
$ws1.Range("B14").Value = "DO NOT ORDER, INCLUDE AND  INSTALL!!!"

# --- Row 25: default row height nudged from 13.5 to 13.8 ---
$ws1.Rows.Item(25).RowHeight = 13.8

# --- New blank spacer rows 27-29 (formatted like D3:F3) merged D27:F29 ---
$ws1.Range("D3:F3").Copy() | Out-Null
$ws1.Range("D27:F27").PasteSpecial(-4122) | Out-Null
$ws1.Range("D28:F28").PasteSpecial(-4122) | Out-Null
$ws1.Range("D29:F29").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(27).RowHeight = 13.8
$ws1.Rows.Item(28).RowHeight = 13.8
$ws1.Rows.Item(29).RowHeight = 13.8
$ws1.Range("D27:F29").Merge() | Out-Null

# --- New rows 30/31 (single framed cell in column C) ---
$ws1.Range("C11").Copy() | Out-Null
$ws1.Range("C30").PasteSpecial(-4122) | Out-Null
$ws1.Range("C31").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(30).RowHeight = 13.8
$ws1.Rows.Item(31).RowHeight = 21

# --- New empty row 32 (height only) ---
$ws1.Rows.Item(32).RowHeight = 13.8

# --- New row 33 (framed row echoing the B11:F11 formatting, left blank) ---
$ws1.Range("B11:F11").Copy() | Out-Null
$ws1.Range("B33:F33").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(33).RowHeight = 13.8

# --- New empty rows 34/35 (height only) ---
$ws1.Rows.Item(34).RowHeight = 13.8
$ws1.Rows.Item(35).RowHeight = 13.8

# --- View changes: zoom to 110% everywhere, Sheet1 selection on G2 ---
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 110

$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 110

$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 110
$ws1.Range("G2").Select() | Out-Null
